# BALP ihe-assuranceLevel StructureDefinition: bump to 1.1.1 and refresh history.
$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
# Version
$metadata.Range("B3").Value = "1.1.1"

# Date
$metadata.Range("B8").Value = "2022-10-21T09:04:31-05:00"

# Description: the SMART Health Cards IG link moved from the build.fhir.org
# staging URL to its published hl7.org/fhir/uv location.
$description = @'
The assuranceLevel element carries various types of Assurance level. May be an Identity Assurance (IAL), an Authentication Assurance Level (AAL), a Federation Assurance Level (FAL), or other. 

In SAML this is [defined to be carried](https://docs.oasis-open.org/security/saml/v2.0/saml-authn-context-2.0-os.pdf) in the `saml:AuthnContextClassRef`, but may be carried elsewhere based on the use-case and profiling of SAML.

The Vocabulary is not defined here. Some sources of vocabulary:
- HL7 v3 [Security Trust Assurance ValueSet](https://terminology.hl7.org/3.0.0/ValueSet-v3-SecurityTrustAssuranceObservationValue.html). These include ISO-7498-2, NIST 800-63-1, and NIST-800-63-2.
- [idmanagement.gov](https://developers.login.gov/saml/#specifying-attributes-and-assurance-levels) published on login.gov
  - this is defined to be carried in the saml:AuthnContextClassRef
- OASIS [Authentication Context for SAML](https://docs.oasis-open.org/security/saml/v2.0/saml-authn-context-2.0-os.pdf)
- An example of a customized and purpose defined small set of codes can be found in the HL7 [SMART Health Cards](http://hl7.org/fhir/uv/shc-vaccination/2021Sep/) defines a [valueset](http://hl7.org/fhir/uv/shc-vaccination/ValueSet/identity-assurance-level)
'@
$metadata.Range("B14").Value = $description

# --- Elements sheet ---------------------------------------------------
# Same description text (Definition column) on the root "Extension" row.
$elements.Range("L2").Value = $description

# The root "Extension" row no longer repeats the ele-1/ext-1 invariant in
# the Constraint(s) column (it still applies via Extension.extension below).
$elements.Range("AI2").Value = ""
